$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the IRON row (row 4) for the new LUCIANO entry
$ws.Rows(4).Insert()
$ws.Cells.Item(4, 1).Value = "'005002390"
$ws.Cells.Item(4, 2).Value = "LUCIANO"
$ws.Cells.Item(4, 3).Value = 11000

# IRON row (now shifted down to row 5): update balance from 9872.65 to 6000
$ws.Cells.Item(5, 3).Value = 6000

# Remove the FLAVIA, ADRIANA, ELAINE, AYRTON, ERIKA rows (now rows 6-10)
$ws.Range("6:10").Delete()

# Remove the FERNANDO row (now row 7, right after GUSTAVO at row 6)
$ws.Rows(7).Delete()
